# daily auto push: 2026-01-27 18:55 UTC
#
# The log sheet gained one new reading for 2026/01/27 at hour 23. It is
# inserted right after the existing 2026/01/27 19:00 entry (old row 725),
# so it becomes the new row 726 and every following row shifts down by one
# (old row 726 -> new 727, ..., old row 767 -> new 768). The sheet's used
# range grows from A1:D767 to A1:D768 accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 726.. down by one to make room for the new reading.
$ws.Rows.Item(726).Insert()

# Column A holds literal text dates ("2026/01/27"), not real Excel dates.
# Force text so Excel's COM layer doesn't auto-convert the date-shaped
# string into a date serial, then restore the plain "Normal" style so no
# stray number-format survives on the cell.
$ws.Range("A726").NumberFormat = "@"
$ws.Range("A726").Value = "2026/01/27"
$ws.Range("A726").Style = "Normal"

$ws.Range("B726").Value = "火"
$ws.Range("C726").Value = 23
$ws.Range("D726").Value = 201
